# Updated cryptos list on Mon Oct  9 03:35:11 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (col D) and "Volume(1h)" (col E) columns for the
# coin rows on the active sheet with newly-scraped values. Both columns
# hold plain text in the workbook (e.g. "210.87", "  -0.62%  "), so every
# write goes through Set-CellText below, which forces a text assignment -
# for values that look like a bare number (e.g. "210.80", "7.60") a
# leading apostrophe is used, exactly like a user typing '210.80 into the
# Excel UI, so Excel doesn't silently reinterpret the text as a Double and
# drop significant trailing zeros / precision.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($Ws, $Cell, $Text, $ForceText) {
    if ($ForceText) {
        $Ws.Range($Cell).Value = "'" + $Text
    } else {
        $Ws.Range($Cell).Value = $Text
    }
}

Set-CellText $ws "D2" '27.810.01' $false
Set-CellText $ws "E2" '  -0.64%  ' $false
Set-CellText $ws "D3" '1.626.54' $false
Set-CellText $ws "E3" '  -0.39%  ' $false
Set-CellText $ws "D4" '0.996' $true
Set-CellText $ws "E4" '  -0.15%  ' $false
Set-CellText $ws "D5" '210.80' $true
Set-CellText $ws "E5" '  -0.61%  ' $false
Set-CellText $ws "D6" '0.517' $true
Set-CellText $ws "E6" '  -0.98%  ' $false
Set-CellText $ws "D7" '0.995' $true
Set-CellText $ws "E7" '  -0.13%  ' $false
Set-CellText $ws "D9" '0.257' $true
Set-CellText $ws "E9" '  -0.61%  ' $false
Set-CellText $ws "E10" '  -1.16%  ' $false
Set-CellText $ws "D11" '0.0879' $true
Set-CellText $ws "E11" '  -0.22%  ' $false
Set-CellText $ws "D12" '1.857.40' $false
Set-CellText $ws "E12" '  -0.38%  ' $false
Set-CellText $ws "D13" '1.631.19' $false
Set-CellText $ws "E13" '  -0.02%  ' $false
Set-CellText $ws "E14" '  -1.29%  ' $false
Set-CellText $ws "E15" '  -1.35%  ' $false
Set-CellText $ws "D16" '64.83' $true
Set-CellText $ws "E16" '  -1.24%  ' $false
Set-CellText $ws "D17" '27.838.44' $false
Set-CellText $ws "E17" '  -0.52%  ' $false
Set-CellText $ws "D18" '227.95' $true
Set-CellText $ws "E18" '  -1.75%  ' $false
Set-CellText $ws "D19" '7.63' $true
Set-CellText $ws "E19" '  +1.07%  ' $false
Set-CellText $ws "E20" '  -1.30%  ' $false
Set-CellText $ws "E21" '  -0.12%  ' $false
Set-CellText $ws "D22" '4.33' $true
Set-CellText $ws "E22" '  -0.63%  ' $false
Set-CellText $ws "D23" '9.92' $true
Set-CellText $ws "E23" '  -5.05%  ' $false
Set-CellText $ws "D24" '2.06' $true
Set-CellText $ws "E24" '  -0.44%  ' $false
Set-CellText $ws "D25" '155.27' $true
Set-CellText $ws "E25" '  +0.60%  ' $false
Set-CellText $ws "E26" '  -0.18%  ' $false
Set-CellText $ws "E27" '  -0.57%  ' $false
Set-CellText $ws "E28" '  -1.32%  ' $false
Set-CellText $ws "E29" '  -0.20%  ' $false
Set-CellText $ws "E30" '  -0.57%  ' $false
Set-CellText $ws "D31" '0.0479' $true
Set-CellText $ws "E31" '  -0.39%  ' $false
Set-CellText $ws "D32" '3.39' $true
Set-CellText $ws "E32" '  -0.22%  ' $false
Set-CellText $ws "E33" '  -0.24%  ' $false
Set-CellText $ws "D34" '1.408.50' $false
Set-CellText $ws "E34" '  +0.08%  ' $false
Set-CellText $ws "D35" '1.61' $true
Set-CellText $ws "E35" '  +2.40%  ' $false
Set-CellText $ws "E36" '  -0.47%  ' $false
Set-CellText $ws "E37" '  -1.34%  ' $false
Set-CellText $ws "E38" '  -1.15%  ' $false
Set-CellText $ws "D39" '0.553' $true
Set-CellText $ws "E39" '  -0.92%  ' $false
Set-CellText $ws "E40" '  -2.44%  ' $false
Set-CellText $ws "E41" '  -0.12%  ' $false
Set-CellText $ws "E42" '  -2.03%  ' $false
Set-CellText $ws "D43" '65.64' $true
Set-CellText $ws "E43" '  -2.02%  ' $false
Set-CellText $ws "E44" '  -0.37%  ' $false
Set-CellText $ws "D45" '5.41' $true
Set-CellText $ws "E45" '  -1.15%  ' $false
Set-CellText $ws "D46" '1.766.64' $false
Set-CellText $ws "E46" '  -0.48%  ' $false
Set-CellText $ws "E47" '  -3.79%  ' $false
Set-CellText $ws "D48" '88.46' $true
Set-CellText $ws "E48" '  +0.35%  ' $false
Set-CellText $ws "E49" '  +0.78%  ' $false
Set-CellText $ws "E50" '  -0.55%  ' $false
Set-CellText $ws "D51" '7.60' $true
Set-CellText $ws "E51" '  +0.68%  ' $false
